$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers E1 (256) and F1 (512) -- stored as text like the
# existing "32"/"64"/"128" headers, matching their bold/bordered style.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Formula = "=""256"""
$ws.Range("F1").Formula = "=""512"""
$ws.Range("E1:F1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4163)

# New data columns E (256) and F (512) for rows 2-65
$data = @(
    @(2, 20981, 83814),
    @(3, 20981, 83814),
    @(4, 14364, 57098),
    @(5, 14364, 57098),
    @(6, 2798, 13066),
    @(7, 2798, 13066),
    @(8, 2948, 13197),
    @(9, 2948, 13197),
    @(10, 28506, 121920),
    @(11, 28506, 121920),
    @(12, 16743, 67373),
    @(13, 16743, 67373),
    @(14, 2813, 13489),
    @(15, 2813, 13489),
    @(16, 2820, 12398),
    @(17, 2820, 12398),
    @(18, 21165, 83808),
    @(19, 21165, 83808),
    @(20, 10442, 41349),
    @(21, 10442, 41349),
    @(22, 3595, 14442),
    @(23, 3595, 14442),
    @(24, 1651, 5672),
    @(25, 1651, 5672),
    @(26, 23293, 95312),
    @(27, 23293, 95312),
    @(28, 10273, 41118),
    @(29, 10273, 41118),
    @(30, 5449, 22623),
    @(31, 5449, 22623),
    @(32, 673, 2001),
    @(33, 673, 2001),
    @(34, 87, 87),
    @(35, 87, 87),
    @(36, 74, 74),
    @(37, 74, 74),
    @(38, 30, 30),
    @(39, 30, 30),
    @(40, 30, 33),
    @(41, 30, 33),
    @(42, 158, 186),
    @(43, 158, 186),
    @(44, 149, 186),
    @(45, 149, 186),
    @(46, 149, 171),
    @(47, 149, 171),
    @(48, 116, 171),
    @(49, 116, 171),
    @(50, 29, 29),
    @(51, 29, 29),
    @(52, 40, 40),
    @(53, 40, 40),
    @(54, 33, 33),
    @(55, 33, 33),
    @(56, 36, 36),
    @(57, 36, 36),
    @(58, 94, 96),
    @(59, 94, 96),
    @(60, 95, 96),
    @(61, 95, 96),
    @(62, 80, 96),
    @(63, 80, 96),
    @(64, 80, 96),
    @(65, 80, 96)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 5).Value = $item[1]
    $ws.Cells.Item($r, 6).Value = $item[2]
}
